# AB#10829 - Fix column sizing for CovidReport table.
# The table originally had an extra empty 975-dxa wide column between the
# "Date" column and the "Test Type" column. Remove that empty column and
# grow the neighbouring column so the overall table width stays the same.
# Also tidy up a couple of variable names ("test_type"/"test_location")
# that had been split across two runs for no reason.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Remove the empty, narrow (975 dxa = 48.75pt) second column and grow the
# column that used to hold "Test Type"/"test_type" so the table keeps its
# overall width (975 + 3997 = 4972 dxa = 248.6pt).
[void]$t.Columns.Item(2).Delete()
$t.Columns.Item(2).Width = 248.6

# Merge the runs that spelled "test_type" and "test_location" across two
# <w:r> elements back into a single run each (two occurrences of each).
[void]$d.Content.Find.Execute("test_type", $true, $false, $false, $false, $false, $true, 1, $false, "test_type", 2)
[void]$d.Content.Find.Execute("test_type", $true, $false, $false, $false, $false, $true, 1, $false, "test_type", 2)

[void]$d.Content.Find.Execute("test_location", $true, $false, $false, $false, $false, $true, 1, $false, "test_location", 2)
[void]$d.Content.Find.Execute("test_location", $true, $false, $false, $false, $false, $true, 1, $false, "test_location", 2)
